$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the base "Temperature [C]" parameter value in R2 from 1500 to 1000.
# Column A (and its shared formulas) reference $R$2, so this ripples
# through every row's cached value automatically on recalculation.
$ws.Range("R2").Value = 1000

# Update the active cell / selection to M16 (was R3).
$ws.Range("M16").Select() | Out-Null
